$d = $word.ActiveDocument

function Rewrite-ParagraphText($para, [string]$newText) {
    # Rebuilds a paragraph's run content from scratch (as a single clean
    # run), discarding any proofErr / stray markup that wrapped the old
    # runs. Keeps the paragraph's own pPr (style / numbering) intact,
    # because InsertParagraphAfter clones the pPr of the paragraph it is
    # called on.
    $para.Range.InsertParagraphAfter()
    $idx = $para.Index
    $newPara = $d.Paragraphs($idx + 1)
    $newPara.Range.Text = $newText
    $para.Range.Delete()
    return $d.Paragraphs($idx)
}

# ---------------------------------------------------------------------
# 5. "Latex/sweave" (was split "Latex/" + spellStart "sweave" spellEnd):
#    collapse to one clean run, then append three new bullet paragraphs
#    (Unit tests / Travis / codecov), moving the _GoBack bookmark to the
#    end of the new last paragraph.
# ---------------------------------------------------------------------
$pLatex = $d.Paragraphs(46)
$pLatex = Rewrite-ParagraphText $pLatex "Latex/sweave"

$pLatex.Range.InsertParagraphAfter()
$pUnit = $d.Paragraphs(47)
$pUnit.Range.Text = "Unit tests"
$pUnit.Range.InsertAfter(" (testthat)")

$pUnit.Range.InsertParagraphAfter()
$pTravis = $d.Paragraphs(48)
$pTravis.Range.Text = "Travis (continuous integration)"

$pTravis.Range.InsertParagraphAfter()
$pCodecov = $d.Paragraphs(49)
$pCodecov.Range.Text = "codecov for telling you which lines of code are tested"

# Move the "_GoBack" bookmark (previously trailing "Meet once that is
# understood") to the very end of the document, after this run.
$d.Bookmarks("\EndOfDoc").Range.Select() | Out-Null
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$d.Bookmarks.Add("_GoBack", $endRange) | Out-Null

# ---------------------------------------------------------------------
# 4. "LaTeX" -- drop the spellStart/spellEnd proofErr wrap.
# ---------------------------------------------------------------------
$pLatex2 = $d.Paragraphs(34)
Rewrite-ParagraphText $pLatex2 "LaTeX" | Out-Null

# ---------------------------------------------------------------------
# 3. "Ypma" -- drop the spellStart/spellEnd proofErr wrap.
# ---------------------------------------------------------------------
$pYpma = $d.Paragraphs(25)
Rewrite-ParagraphText $pYpma "Ypma" | Out-Null

# ---------------------------------------------------------------------
# 2. "Meet once that is understood" -- remove the old _GoBack bookmark
#    that used to sit here (it has been re-created at the new end of
#    the document above).
# ---------------------------------------------------------------------
$pMeet = $d.Paragraphs(23)
Rewrite-ParagraphText $pMeet "Meet once that is understood" | Out-Null

# ---------------------------------------------------------------------
# 1. "git scm for documentation" -- merge the three runs (with the
#    spellStart/spellEnd wrapped "scm" run in the middle) into one.
# ---------------------------------------------------------------------
$pGit = $d.Paragraphs(14)
Rewrite-ParagraphText $pGit "git scm for documentation" | Out-Null

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
